# Updated cryptos list on Mon Mar  4 05:03:43 UTC 2024 with GitHub Actions
# Refresh the Price / Volume(1h) columns for each coin row, and re-sort a
# few rows whose relative ranking changed (FirstDigitalUSD/WEMIXToken swap,
# and the NEARProtocol/LidoDAOToken/Monero/PEPE block shifting down as PEPE
# jumped to the top of that group).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat "@" (Text) is applied before writing any cell whose new
# value looks like a plain number (e.g. "415.13"), so Excel keeps storing
# it as text instead of silently coercing it to a numeric cell.
$ws.Range("D2").Value = '63.636.73'
$ws.Range("E2").Value = '  +2.55%  '
$ws.Range("D3").Value = '3.484.67'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.13'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.70'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("E7").Value = '  -0.77%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.745'
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.157'
$ws.Range("E10").Value = '  +10.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.47'
$ws.Range("E11").Value = '  -3.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.74'
$ws.Range("E12").Value = '  +4.00%  '
$ws.Range("E13").Value = '  +2.67%  '
$ws.Range("D14").Value = '4.040.76'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.31'
$ws.Range("E16").Value = '  -4.77%  '
$ws.Range("D17").Value = '3.499.94'
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.44'
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").Value = '63.570.52'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '458.29'
$ws.Range("E21").Value = '  -7.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.22'
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.30'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.28'
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.19'
$ws.Range("E25").Value = '  +10.45%  '
$ws.Range("E26").Value = '  -3.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.55'
$ws.Range("E27").Value = '  -4.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.76'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.49'
$ws.Range("E29").Value = '  +2.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.51'
$ws.Range("E30").Value = '  -1.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.69'
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.13'
$ws.Range("E34").Value = '  -4.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.27'
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0488'
$ws.Range("E37").Value = '  -2.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  +4.02%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("E39").Value = '  +4.76%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0653'
$ws.Range("E42").Value = '  +53.38%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.51'
$ws.Range("E43").Value = '  +3.94%  '
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.35'
$ws.Range("E44").Value = '  -4.39%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.90'
$ws.Range("E45").Value = '  -3.60%  '
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("E47").Value = '  -7.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.35'
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.07'
$ws.Range("E49").Value = '  -3.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.72'
$ws.Range("E50").Value = '  -5.89%  '
$ws.Range("E51").Value = '  -5.85%  '

Write-Host "Applied cryptos update"